# Daily attendance processing - 2026-01-02 19:27:28
# Swap the order of names in column G ("Recorded By") from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value2 = $newVal
    }
}
